# Update Sheets via scheduled runner: refresh market-price derived columns (H-N)
# for the affected Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 751.6667
$ws.Range("I6").Value = 758
$ws.Range("J6").Value = 720
$ws.Range("K6").Value = 2274
$ws.Range("L6").Value = 2160
$ws.Range("M6").Value = -2162
$ws.Range("N6").Value = -2384
# Row 80
$ws.Range("H80").Value = 908.5
$ws.Range("I80").Value = 1048
$ws.Range("J80").Value = 848.7143
$ws.Range("K80").Value = 3144
$ws.Range("L80").Value = 2546.1429
$ws.Range("M80").Value = -2146
$ws.Range("N80").Value = -4542.1429
# Row 83
$ws.Range("H83").Value = 908.5
$ws.Range("I83").Value = 1048
$ws.Range("J83").Value = 848.7143
$ws.Range("K83").Value = 9432
$ws.Range("L83").Value = 7638.428699999999
$ws.Range("M83").Value = -4440
$ws.Range("N83").Value = -17622.4287
# Row 86
$ws.Range("H86").Value = 2945.9524
$ws.Range("I86").Value = 1733.2354
$ws.Range("J86").Value = 8100
$ws.Range("K86").Value = 1733.2354
$ws.Range("L86").Value = 8100
$ws.Range("M86").Value = -610.2354
# Row 89
$ws.Range("H89").Value = 2945.9524
$ws.Range("I89").Value = 1733.2354
$ws.Range("J89").Value = 8100
$ws.Range("K89").Value = 8666.177
$ws.Range("L89").Value = 40500
$ws.Range("M89").Value = -3050.177
# Row 111
$ws.Range("H111").Value = 8466956
$ws.Range("I111").Value = 11906198
$ws.Range("J111").Value = 212773.6
$ws.Range("K111").Value = 35718594
$ws.Range("L111").Value = 638320.8
$ws.Range("M111").Value = -35715527
$ws.Range("N111").Value = -644454.8
# Row 112
$ws.Range("H112").Value = 3428.6287
$ws.Range("I112").Value = 1580.5
$ws.Range("J112").Value = 3811
$ws.Range("K112").Value = 4741.5
$ws.Range("L112").Value = 11433
$ws.Range("M112").Value = -3633.5
$ws.Range("N112").Value = -13649
# Row 135
$ws.Range("H135").Value = 1104.4138
$ws.Range("I135").Value = 926.96295
$ws.Range("J135").Value = 3500
$ws.Range("K135").Value = 8342.66655
$ws.Range("L135").Value = 31500
$ws.Range("M135").Value = -5807.66655
# Row 137
$ws.Range("H137").Value = 2597.7932
$ws.Range("I137").Value = 2178
$ws.Range("J137").Value = 3699.75
$ws.Range("K137").Value = 6534
$ws.Range("L137").Value = 11099.25
$ws.Range("M137").Value = -3984
$ws.Range("N137").Value = -16199.25
# Row 138
$ws.Range("H138").Value = 2500.1458
$ws.Range("I138").Value = 1354.05
$ws.Range("J138").Value = 3318.7856
$ws.Range("K138").Value = 4062.15
$ws.Range("L138").Value = 9956.356800000001
$ws.Range("M138").Value = 1077.85
$ws.Range("N138").Value = -20236.3568

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 416.35715
$ws.Range("I26").Value = 744.75
$ws.Range("J26").Value = 285
$ws.Range("K26").Value = 744.75
$ws.Range("L26").Value = 285
$ws.Range("M26").Value = -414.75
$ws.Range("N26").Value = -945
# Row 45
$ws.Range("H45").Value = 100002430
$ws.Range("I45").Value = 137502020
$ws.Range("J45").Value = 3569
$ws.Range("K45").Value = 137502020
$ws.Range("L45").Value = 3569
$ws.Range("M45").Value = -137501643
# Row 74
$ws.Range("H74").Value = 4722.077
$ws.Range("I74").Value = 2641.647
$ws.Range("J74").Value = 8651.777
$ws.Range("K74").Value = 2641.647
$ws.Range("L74").Value = 8651.777
$ws.Range("M74").Value = -1767.647
# Row 77
$ws.Range("H77").Value = 4722.077
$ws.Range("I77").Value = 2641.647
$ws.Range("J77").Value = 8651.777
$ws.Range("K77").Value = 13208.235
$ws.Range("L77").Value = 43258.885
$ws.Range("M77").Value = -8840.235000000001
# Row 122
$ws.Range("H122").Value = 2684.6978
$ws.Range("I122").Value = 2572.2285
$ws.Range("J122").Value = 3176.75
$ws.Range("K122").Value = 7716.685500000001
$ws.Range("L122").Value = 9530.25
$ws.Range("M122").Value = -5266.685500000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2618.7144
$ws.Range("I107").Value = 2619.8
$ws.Range("J107").Value = 2618.111
$ws.Range("K107").Value = 2619.8
$ws.Range("L107").Value = 2618.111
$ws.Range("M107").Value = -699.8000000000002
$ws.Range("N107").Value = -6458.111
# Row 134
$ws.Range("H134").Value = 2550.907
$ws.Range("I134").Value = 1986.8948
$ws.Range("J134").Value = 6837.4
$ws.Range("K134").Value = 5960.6844
$ws.Range("L134").Value = 20512.2
$ws.Range("M134").Value = -3425.6844

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 2385.6365
$ws.Range("I99").Value = 3135.818
$ws.Range("J99").Value = 1635.4546
$ws.Range("K99").Value = 3135.818
$ws.Range("L99").Value = 1635.4546
$ws.Range("M99").Value = -1637.818
# Row 105
$ws.Range("H105").Value = 1379
$ws.Range("I105").Value = 1661.8334
$ws.Range("J105").Value = 530.5
$ws.Range("K105").Value = 1661.8334
$ws.Range("L105").Value = 530.5
$ws.Range("M105").Value = 85.16660000000002
$ws.Range("N105").Value = -4024.5
# Row 107
$ws.Range("H107").Value = 6946003.5
$ws.Range("I107").Value = 1397.3334
$ws.Range("J107").Value = 18520348
$ws.Range("K107").Value = 1397.3334
$ws.Range("L107").Value = 18520348
$ws.Range("M107").Value = 522.6666
$ws.Range("N107").Value = -18524188
# Row 126
$ws.Range("H126").Value = 2385.6365
$ws.Range("I126").Value = 3135.818
$ws.Range("J126").Value = 1635.4546
$ws.Range("K126").Value = 9407.454000000002
$ws.Range("L126").Value = 4906.3638
$ws.Range("M126").Value = -6937.454000000002
# Row 134
$ws.Range("H134").Value = 3520.5
$ws.Range("I134").Value = 2723.25
$ws.Range("J134").Value = 5912.25
$ws.Range("K134").Value = 8169.75
$ws.Range("L134").Value = 17736.75
$ws.Range("M134").Value = -5634.75
$ws.Range("N134").Value = -22806.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1260.5834
$ws.Range("I113").Value = 1294
$ws.Range("J113").Value = 1253.9
$ws.Range("K113").Value = 3882
$ws.Range("L113").Value = 3761.7
$ws.Range("M113").Value = -1712
$ws.Range("N113").Value = -8101.700000000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 80
$ws.Range("H80").Value = 2608
$ws.Range("I80").Value = 2760
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 2760
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -1762
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 83
$ws.Range("H83").Value = 2608
$ws.Range("I83").Value = 2760
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 13800
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -8808
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 107
$ws.Range("H107").Value = 16667188
$ws.Range("I107").Value = 27778014
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 27778014
$ws.Range("L107").Value = 950
$ws.Range("M107").Value = -27776094
$ws.Range("N107").Value = -4790
# Row 126
$ws.Range("H126").Value = 3414.9512
$ws.Range("I126").Value = 3303.1316
$ws.Range("J126").Value = 4831.3335
$ws.Range("K126").Value = 9909.3948
$ws.Range("L126").Value = 14494.0005
$ws.Range("M126").Value = -7439.3948
# Row 132
$ws.Range("H132").Value = 3951.678
$ws.Range("I132").Value = 3698.8235
$ws.Range("J132").Value = 4295.56
$ws.Range("K132").Value = 11096.4705
$ws.Range("L132").Value = 12886.68
$ws.Range("M132").Value = -8566.470499999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6540.6665
$ws.Range("I40").Value = 6193.222
$ws.Range("J40").Value = 7583
$ws.Range("K40").Value = 6193.222
$ws.Range("L40").Value = 7583
$ws.Range("M40").Value = -6057.222
# Row 122
$ws.Range("H122").Value = 6369.643
$ws.Range("I122").Value = 5567.7
$ws.Range("J122").Value = 8374.5
$ws.Range("K122").Value = 16703.1
$ws.Range("L122").Value = 25123.5
$ws.Range("M122").Value = -14253.1

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 3015
$ws.Range("I32").Value = 3015
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3015
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2698
$ws.Range("N32").ClearContents()
# Row 93
$ws.Range("H93").Value = 79999
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 79999
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 79999
$ws.Range("N93").Value = -84991
# Row 136
$ws.Range("H136").Value = 5630.6
$ws.Range("I136").Value = 4679.091
$ws.Range("J136").Value = 8247.25
$ws.Range("K136").Value = 14037.273
$ws.Range("L136").Value = 24741.75
$ws.Range("M136").Value = -11487.273
